$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Unsuccessful login scenario
$ws.Range("A3").Value = "Unsuccessfull login with invalid credentials"
$ws.Range("B3").Value = "testingInvald"
$ws.Range("C3").Value = "cvjdvjjvj"

# Row 4: Login with empty fields (only first column populated)
$ws.Range("A4").Value = "Login with empty fields"

# Row 5: User log out scenario (reuses existing "Admin" / "admin123" strings)
$ws.Range("A5").Value = "User log out Successfully"
$ws.Range("B5").Value = "Admin"
$ws.Range("C5").Value = "admin123"

# Update selection to match the final active cell used in the source workbook
[void]$ws.Range("C5").Select()
